$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A2 value change: "r" -> "sr"
$ws.Range("A2").Value = "sr"

# New columns E/F for rows 4-14 (numbering 0..10 plus new labels)
$labels1 = @("sub det","sub noun","relation","obj det","obj noun","space relat","space det","space obj","time relat","time det","time obj")
for ($i = 0; $i -lt 11; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 5).Value = $i
    $ws.Cells.Item($row, 6).Value = $labels1[$i]
}

# New rows 15-20 in columns A/B/C
$ws.Range("A15").Value = "n"
$ws.Range("B15").Value = "time"

$ws.Range("A16").Value = "n"
$ws.Range("B16").Value = "moment"

$ws.Range("A17").Value = "r"
$ws.Range("B17").Value = "moves"
$ws.Range("C17").Value = "MV"

$ws.Range("A18").Value = "tr"
$ws.Range("B18").Value = "through"
$ws.Range("C18").Value = "TRG"

$ws.Range("A19").Value = "d"
$ws.Range("B19").Value = "a"

$ws.Range("A20").Value = "d"
$ws.Range("B20").Value = "the"

# New columns E/F for rows 15-18 (numbering 11..14 plus new "changed" labels)
$labels2 = @("sub changed","obj changed","space obj changed","time obj changed")
for ($i = 0; $i -lt 4; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 5).Value = 11 + $i
    $ws.Cells.Item($row, 6).Value = $labels2[$i]
}

# Page setup (portrait orientation, matching OOXML pageSetup addition)
$ws.PageSetup.Orientation = 1

# Move the active selection to the last-entered cell
$ws.Range("E19").Select()
